# Build site at 2022-09-26 16:07:08 UTC
# Restructure rows 13-24 of the LOM3089 syllabus sheet:
#  - row 13 gains a label in A ("Programa resumido:") and its B/C value
#    becomes "Semestral" (was "519033 - Carlos Yujiro Shigue")
#  - several long essay-style texts (Objetivos longo, Programa resumido
#    longo, Programa longo, Bibliografia longa) are removed
#  - the remaining label/value pairs shift up to fill the gaps
#  - the very last row (24) disappears, shrinking the sheet to C23

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: "Programa resumido:" / "Semestral" --------------------------
$ws.Cells.Item(13, 1).Value = "Programa resumido:"
$ws.Cells.Item(13, 2).Value = "Semestral"
$ws.Cells.Item(13, 3).Value = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# --- Row 14: "Short syllabus:" (B/C cleared) ------------------------------
$ws.Cells.Item(14, 1).Value = "Short syllabus:"
$ws.Cells.Item(14, 2).ClearContents()
$ws.Cells.Item(14, 3).ClearContents()
$ws.Rows.Item(14).RowHeight = 60

# --- Row 15: "Programa:" / "01/01/2012" -----------------------------------
$ws.Cells.Item(15, 1).Value = "Programa:"
$ws.Cells.Item(15, 2).Value = "01/01/2012"
$ws.Cells.Item(15, 3).Value = "01/01/2012"
$ws.Rows.Item(15).RowHeight = 120

# --- Row 16: "Syllabus:" (B/C cleared) ------------------------------------
$ws.Cells.Item(16, 1).Value = "Syllabus:"
$ws.Cells.Item(16, 2).ClearContents()
$ws.Cells.Item(16, 3).ClearContents()
$ws.Rows.Item(16).RowHeight = 120

# --- Row 17: "Avaliação:" only, back to default (non-custom) height ------
$ws.Cells.Item(17, 1).Value = "Avaliação:"
$ws.Cells.Item(17, 2).ClearContents()
$ws.Cells.Item(17, 3).ClearContents()
$ws.Rows.Item(17).EntireRow.AutoFit()

# --- Row 18: "Método:" / "519033 - Carlos Yujiro Shigue" ------------------
$ws.Cells.Item(18, 1).Value = "Método:"
$ws.Cells.Item(18, 2).Value = "519033 - Carlos Yujiro Shigue"
$ws.Cells.Item(18, 3).Value = "519033 - Carlos Yujiro Shigue"
$ws.Rows.Item(18).RowHeight = 60

# --- Row 19: "Critério:" / evaluation method text -------------------------
$ws.Cells.Item(19, 1).Value = "Critério:"
$ws.Cells.Item(19, 2).Value = "A avaliação será feita por meio de duas provas escritas P1 e P2 e por listas de exercícios e relatórios."
$ws.Cells.Item(19, 3).Value = "A avaliação será feita por meio de duas provas escritas P1 e P2 e por listas de exercícios e relatórios."
$ws.Rows.Item(19).RowHeight = 60

# --- Row 20: "Norma de recuperação:" / grade formula text -----------------
$ws.Cells.Item(20, 1).Value = "Norma de recuperação:"
$ws.Cells.Item(20, 2).Value = "A Nota final (NF) será calculada pela média ponderada das provas escritas e pela média dos trabalhos TR da seguinte maneira: NF = (P1 + 2*P2 + TR)/4"
$ws.Cells.Item(20, 3).Value = "A Nota final (NF) será calculada pela média ponderada das provas escritas e pela média dos trabalhos TR da seguinte maneira: NF = (P1 + 2*P2 + TR)/4"
$ws.Rows.Item(20).RowHeight = 60

# --- Row 21: "Bibliografia:" / recovery formula text -----------------------
$ws.Cells.Item(21, 1).Value = "Bibliografia:"
$ws.Cells.Item(21, 2).Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"
$ws.Cells.Item(21, 3).Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"
$ws.Rows.Item(21).RowHeight = 120

# --- Row 22: "Requisitos:" only, back to default height -------------------
$ws.Cells.Item(22, 1).Value = "Requisitos:"
$ws.Cells.Item(22, 2).ClearContents()
$ws.Cells.Item(22, 3).ClearContents()
$ws.Rows.Item(22).EntireRow.AutoFit()

# --- Row 23: requirement detail moves up from row 24, A cleared -----------
$ws.Cells.Item(23, 1).ClearContents()
$ws.Cells.Item(23, 2).Value = "LOB1019 -  Física II  (Requisito fraco)`n"
$ws.Cells.Item(23, 3).Value = "LOB1019 -  Física II  (Requisito fraco)`n"
$ws.Rows.Item(23).RowHeight = 30

# --- Row 24 no longer exists: delete it so dimension shrinks to C23 -------
$ws.Rows.Item(24).Delete()
